# Add season-record columns (Wins, Losses, Ties) to the team/player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (columns AC, AD, AE)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) used by A1:AB1
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill in the season record (91 wins, 71 losses, 0 ties) for every data row
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 29).Value = 91
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 0
}
